# Auto-generated script applying value updates per the OOXML diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 0
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 0
$ws.Range("L16").ClearContents()
$ws.Range("M16").ClearContents()
$ws.Range("N16").Value = 0

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H49").Value = 0
$ws.Range("I49").Value = 0
$ws.Range("K49").Value = 0
$ws.Range("M49").ClearContents()

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H58").Value = 184
$ws.Range("I58").Value = 189.22223
$ws.Range("J58").Value = 137
$ws.Range("K58").Value = 567.66669
$ws.Range("L58").Value = 411
$ws.Range("M58").Value = -417.66669
$ws.Range("N58").Value = -711

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H69").Value = 9750
$ws.Range("J69").Value = 9750
$ws.Range("L69").Value = 29250
$ws.Range("N69").Value = -30998

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H72").Value = 9750
$ws.Range("J72").Value = 9750
$ws.Range("L72").Value = 87750
$ws.Range("N72").Value = -96486

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 5372
$ws.Range("I116").Value = 5119.6
$ws.Range("K116").Value = 5119.6
$ws.Range("M116").Value = -1677.6

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7458.606
$ws.Range("I32").Value = 7458.606
$ws.Range("K32").Value = 7458.606
$ws.Range("M32").Value = -7171.606

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H80").Value = 37994.9
$ws.Range("J80").Value = 37995
$ws.Range("L80").Value = 37995
$ws.Range("N80").Value = -39991

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H83").Value = 37994.9
$ws.Range("J83").Value = 37995
$ws.Range("L83").Value = 113985
$ws.Range("N83").Value = -123969

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 1699.6666
$ws.Range("I122").Value = 1549.5
$ws.Range("J122").Value = 2000
$ws.Range("K122").Value = 4648.5
$ws.Range("L122").Value = 6000
$ws.Range("M122").Value = -2198.5
$ws.Range("N122").Value = -10900

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 7999.5
$ws.Range("I132").Value = 7999.5
$ws.Range("K132").Value = 23998.5
$ws.Range("M132").Value = -21468.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H11").Value = 790.6
$ws.Range("I11").Value = 502
$ws.Range("J11").Value = 983
$ws.Range("K11").Value = 502
$ws.Range("L11").Value = 983
$ws.Range("M11").Value = -362
$ws.Range("N11").Value = -1263

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1127.6666
$ws.Range("I86").Value = 1156.125
$ws.Range("K86").Value = 1156.125
$ws.Range("M86").Value = -33.125

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 1127.6666
$ws.Range("I89").Value = 1156.125
$ws.Range("K89").Value = 5780.625
$ws.Range("M89").Value = -164.625

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H122").Value = 36999
$ws.Range("J122").Value = 36999
$ws.Range("L122").Value = 36999
$ws.Range("N122").Value = -46799

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H123").Value = 49999
$ws.Range("J123").Value = 49999
$ws.Range("L123").Value = 49999
$ws.Range("N123").Value = -59799

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 20000
$ws.Range("I134").Value = 20000
$ws.Range("K134").Value = 60000
$ws.Range("M134").Value = -57465

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H33").Value = 2760.1667
$ws.Range("I33").Value = 2612.2
$ws.Range("J33").Value = 3500
$ws.Range("K33").Value = 2612.2
$ws.Range("L33").Value = 3500
$ws.Range("M33").Value = -2233.2
$ws.Range("N33").Value = -4258

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("L51").ClearContents()
$ws.Range("N51").Value = 0

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H55").Value = 8500
$ws.Range("I55").Value = 7750
$ws.Range("K55").Value = 7750
$ws.Range("M55").Value = -7435

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2002.75
$ws.Range("J58").Value = 1999.5
$ws.Range("L58").Value = 1999.5
$ws.Range("N58").Value = -2405.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H61").Value = 0
$ws.Range("J61").Value = 0
$ws.Range("L61").ClearContents()
$ws.Range("N61").Value = 0

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H74").Value = 0
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 0
$ws.Range("L74").ClearContents()
$ws.Range("M74").ClearContents()
$ws.Range("N74").Value = 0

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H77").Value = 0
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 0
$ws.Range("L77").ClearContents()
$ws.Range("M77").ClearContents()
$ws.Range("N77").Value = 0

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 2000
$ws.Range("I99").Value = 2000
$ws.Range("K99").Value = 2000
$ws.Range("M99").Value = -502

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 2000
$ws.Range("I126").Value = 2000
$ws.Range("K126").Value = 6000
$ws.Range("M126").Value = -3530

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 2935.6428
$ws.Range("I134").Value = 2766.5833
$ws.Range("K134").Value = 8299.749899999999
$ws.Range("M134").Value = -5764.749899999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 2002.75
$ws.Range("J136").Value = 1999.5
$ws.Range("L136").Value = 5998.5
$ws.Range("N136").Value = -11098.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 3790887.5
$ws.Range("I8").Value = 3790887.5
$ws.Range("K8").Value = 11372662.5
$ws.Range("M8").Value = -11372523.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H97").Value = 148.66667
$ws.Range("I97").Value = 173.5
$ws.Range("J97").Value = 99
$ws.Range("K97").Value = 520.5
$ws.Range("L97").Value = 297
$ws.Range("M97").Value = -24.5
$ws.Range("N97").Value = -1289

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H117").Value = 1681.2
$ws.Range("I117").Value = 1844.5
$ws.Range("K117").Value = 5533.5
$ws.Range("M117").Value = -2091.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 168831.67
$ws.Range("I3").Value = 168831.67
$ws.Range("K3").Value = 168831.67
$ws.Range("M3").Value = -168715.67

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("M23").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3039.4285
$ws.Range("I102").Value = 2504.818
$ws.Range("J102").Value = 4999.6665
$ws.Range("K102").Value = 2504.818
$ws.Range("L102").Value = 4999.6665
$ws.Range("M102").Value = -882.8180000000002
$ws.Range("N102").Value = -8243.666499999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 4783.25
$ws.Range("J132").Value = 7500
$ws.Range("L132").Value = 22500
$ws.Range("N132").Value = -27560

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1354.3334
$ws.Range("J16").Value = 2999.5
$ws.Range("L16").Value = 2999.5
$ws.Range("N16").Value = -3339.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H35").Value = 3000
$ws.Range("I35").Value = 3000
$ws.Range("K35").Value = 3000
$ws.Range("M35").Value = -2664

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 1488.5
$ws.Range("I93").Value = 1541.8
$ws.Range("J93").Value = 1399.6666
$ws.Range("K93").Value = 1541.8
$ws.Range("L93").Value = 1399.6666
$ws.Range("M93").Value = -293.8
$ws.Range("N93").Value = -3895.6666

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 6828.8335
$ws.Range("I100").Value = 6828.8335
$ws.Range("K100").Value = 6828.8335
$ws.Range("M100").Value = -6287.8335

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 3850.9092
$ws.Range("I122").Value = 3850.9092
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 11552.7276
$ws.Range("L122").Value = 0
$ws.Range("M122").ClearContents()
$ws.Range("N122").Value = -9102.7276

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 2749.25
$ws.Range("J4").Value = 2749.25
$ws.Range("L4").Value = 2749.25
$ws.Range("N4").Value = -2975.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 6128.5386
$ws.Range("I107").Value = 5355.875
$ws.Range("K107").Value = 16067.625
$ws.Range("M107").Value = -14147.625

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 7001.75
$ws.Range("J132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 4424.15
$ws.Range("I136").Value = 4424.15
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 13272.45
$ws.Range("L136").Value = 0
$ws.Range("M136").ClearContents()
$ws.Range("N136").Value = -10722.45
